{"js": "// Replace the two-digit-multiplication equations with their new values.\n// Each \"before\" string is unique in the document, so a plain text search\n// (no wildcards) for each one followed by a full-text replace is safe.\nconst replacements = [\n  [\"41\u00d730=1230\", \"43\u00d790=3870\"],\n  [\"77\u00d799=7623\", \"43\u00d714=602\"],\n  [\"94\u00d715=1410\", \"80\u00d756=4480\"],\n  [\"19\u00d760=1140\", \"68\u00d765=4420\"],\n  [\"98\u00d775=7350\", \"99\u00d721=2079\"],\n  [\"71\u00d785=6035\", \"27\u00d751=1377\"],\n  [\"45\u00d781=3645\", \"41\u00d764=2624\"],\n  [\"62\u00d778=4836\", \"72\u00d757=4104\"],\n  [\"48\u00d752=2496\", \"74\u00d716=1184\"],\n  [\"50\u00d764=3200\", \"53\u00d739=2067\"],\n  [\"51\u00d731=1581\", \"65\u00d793=6045\"],\n  [\"34\u00d753=1802\", \"22\u00d770=1540\"],\n  [\"83\u00d715=1245\", \"33\u00d746=1518\"],\n  [\"81\u00d752=4212\", \"18\u00d746=828\"],\n  [\"78\u00d772=5616\", \"43\u00d782=3526\"],\n  [\"16\u00d798=1568\", \"14\u00d717=238\"],\n  [\"14\u00d758=812\", \"55\u00d720=1100\"],\n  [\"80\u00d730=2400\", \"79\u00d779=6241\"],\n  [\"73\u00d735=2555\", \"27\u00d711=297\"],\n  [\"75\u00d732=2400\", \"61\u00d774=4514\"],\n  [\"98\u00d757=5586\", \"44\u00d753=2332\"],\n  [\"99\u00d746=4554\", \"48\u00d768=3264\"],\n  [\"52\u00d759=3068\", \"59\u00d773=4307\"],\n  [\"74\u00d748=3552\", \"59\u00d775=4425\"],\n  [\"32\u00d732=1024\", \"89\u00d768=6052\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the two-digit-multiplication equations with their new values.\n# Each \"before\" string is unique in the document, so Find/Replace targeting\n# each equation individually (wdReplaceAll, but each only matches once) is safe.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"41\u00d730=1230\", \"43\u00d790=3870\"),\n    @(\"77\u00d799=7623\", \"43\u00d714=602\"),\n    @(\"94\u00d715=1410\", \"80\u00d756=4480\"),\n    @(\"19\u00d760=1140\", \"68\u00d765=4420\"),\n    @(\"98\u00d775=7350\", \"99\u00d721=2079\"),\n    @(\"71\u00d785=6035\", \"27\u00d751=1377\"),\n    @(\"45\u00d781=3645\", \"41\u00d764=2624\"),\n    @(\"62\u00d778=4836\", \"72\u00d757=4104\"),\n    @(\"48\u00d752=2496\", \"74\u00d716=1184\"),\n    @(\"50\u00d764=3200\", \"53\u00d739=2067\"),\n    @(\"51\u00d731=1581\", \"65\u00d793=6045\"),\n    @(\"34\u00d753=1802\", \"22\u00d770=1540\"),\n    @(\"83\u00d715=1245\", \"33\u00d746=1518\"),\n    @(\"81\u00d752=4212\", \"18\u00d746=828\"),\n    @(\"78\u00d772=5616\", \"43\u00d782=3526\"),\n    @(\"16\u00d798=1568\", \"14\u00d717=238\"),\n    @(\"14\u00d758=812\", \"55\u00d720=1100\"),\n    @(\"80\u00d730=2400\", \"79\u00d779=6241\"),\n    @(\"73\u00d735=2555\", \"27\u00d711=297\"),\n    @(\"75\u00d732=2400\", \"61\u00d774=4514\"),\n    @(\"98\u00d757=5586\", \"44\u00d753=2332\"),\n    @(\"99\u00d746=4554\", \"48\u00d768=3264\"),\n    @(\"52\u00d759=3068\", \"59\u00d773=4307\"),\n    @(\"74\u00d748=3552\", \"59\u00d775=4425\"),\n    @(\"32\u00d732=1024\", \"89\u00d768=6052\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n\n$d.Save()\n"}
